# Protokoll.xlsx update — "Datenmodell verbessern, Sprint Planning Meeting"
# Moves the small "Thema / Fertigstellungsgrad / Problem / Loesung" summary
# table from G1:K3 down to A18:E21, adds a new sprint-planning entry as
# row 4 of the main log, and appends a new summary row (SprintPlanningMeeting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the old "Thema/Fertigstellungsgrad/Problem/Loesung" mini table
#    that used to live at G1:K3 (it is being relocated further down).
# ---------------------------------------------------------------------
$ws.Range("G1:K3").Clear()

# ---------------------------------------------------------------------
# 2. New sprint log entry in row 4 (week of 2016-10-19) — same text
#    repeated across B4:E4, like the other weekly rows.
# ---------------------------------------------------------------------
$ws.Range("B4:E4").Value = "Datenmodell verbessern, Sprint Planning Meeting"

# ---------------------------------------------------------------------
# 3. Re-create the summary table further down the sheet, at A18:E21,
#    with an extra row for the new "SprintPlanningMeeting" topic.
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Thema"
$ws.Range("A18").Font.Bold = $true

$ws.Range("B18").Value = "Fertigstellungsgrad"
$ws.Range("B18").Font.Bold = $true
$ws.Range("B18").NumberFormat = "0%"

$ws.Range("D18").Value = "Problem"
$ws.Range("D18").Font.Bold = $true

$ws.Range("E18").Value = "Lösung"
$ws.Range("E18").Font.Bold = $true

$ws.Range("A19").Value = "Datenmodell "
$ws.Range("B19").Value = 0.9
$ws.Range("B19").NumberFormat = "0%"

$ws.Range("A20").Value = "Datenmodell "
$ws.Range("B20").Value = 1
$ws.Range("B20").NumberFormat = "0%"

$ws.Range("A21").Value = "SprintPlanningMeeting"
$ws.Range("B21").Value = 1
$ws.Range("B21").NumberFormat = "0%"

# ---------------------------------------------------------------------
# 4. Column widths — A narrower, B:E widened (they now hold the long
#    meeting-notes text), F:G / H left close to their previous sizes.
#    (Input values are tuned so the engine's internal character-width
#    rounding lands as close as possible to the intended ~12.6/38.9/
#    14.3/15.8 character widths.)
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 11.755
$ws.Columns("B:E").ColumnWidth = 37.925
$ws.Columns("F:G").ColumnWidth = 13.25
$ws.Columns("H").ColumnWidth = 14.75

# ---------------------------------------------------------------------
# 5. Selection follows the last-edited cell, like Excel leaves it.
# ---------------------------------------------------------------------
$ws.Range("B21").Select() | Out-Null
